$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Sara Isabel Pisoni"
$ws.Range("B8").Value = "Elia Battisti | U.SGUARNA"
$ws.Range("C8").Value = "Filippo Benetti | I Magnifici"
$ws.Range("D8").Value = "Andrea Conzatti | FC SAVIGNANO"
$ws.Range("E8").Value = "Lorenzo Rossi | Power Ginger"
$ws.Range("F8").Value = "Andreas Galli | SdrumALA"
